# Daily attendance processing - 2026-01-21 16:59:42
# Normalizes the ordering of the comma-separated "Recorded By" list in
# column G of the Session Analysis Results sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$map = @{
    "System, backup@backdoor.com, system" = "system, System, backup@backdoor.com"
    "System, dnasr281@gmail.com"          = "dnasr281@gmail.com, System"
    "System, admin@admin.com"             = "admin@admin.com, System"
    "dnasr281@gmail.com, System"          = "System, dnasr281@gmail.com"
    "dnasr281@gmail.com, admin@admin.com" = "admin@admin.com, dnasr281@gmail.com"
}

for ($r = 2; $r -le 157; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    $cur = $cell.Text
    if ($map.ContainsKey($cur)) {
        $cell.Value = $map[$cur]
    }
}
